$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.237.15"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -0.49%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.862.33"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -0.61%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.000"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.04%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "236.57"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +0.33%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4682"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.25%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2897"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +2.01%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06539"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -0.16%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "21.62"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +2.77%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07937"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +0.10%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "97.84"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +0.56%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.869.60"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -0.24%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.163"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +0.22%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6806"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +0.89%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "267.55"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -5.17%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "30.219.26"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -0.56%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "13.78"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +8.53%  "

$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +0.02%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000007395"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +1.40%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "2.111.68"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -0.73%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.314"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -4.41%  "

$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -0.03%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.175"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -0.50%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "167.02"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +1.75%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.217"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -0.82%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.87"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -1.18%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.958"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +1.11%  "

$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +2.42%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.09835"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +1.53%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.373"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -1.37%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.472"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -0.24%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.045"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -1.77%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.04708"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +0.02%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.131"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +1.08%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7031"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -0.26%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.709"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -0.40%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01875"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +0.94%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.615"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +2.71%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.277"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -1.25%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "74.36"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +0.82%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.936"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -0.44%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.8466"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -0.07%  "

$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +0.00%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.4161"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -0.66%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "103.09"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -0.60%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "958.87"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +3.05%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.150"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -0.96%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.168"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -0.87%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "34.12"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -0.04%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05654"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +0.52%  "
